$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 9199.25
$ws.Range("I18").Value = 11965.667
$ws.Range("J18").Value = 900
$ws.Range("K18").Value = 11965.667
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = -11681.667
$ws.Range("N18").Value = -1468
$ws.Range("H62").Value = 90943680
$ws.Range("I62").Value = 200002990
$ws.Range("K62").Value = 200002990
$ws.Range("M62").Value = -200002366
$ws.Range("H65").Value = 90943680
$ws.Range("I65").Value = 200002990
$ws.Range("K65").Value = 1000014950
$ws.Range("M65").Value = -1000011830
$ws.Range("H86").Value = 118059910
$ws.Range("I86").Value = 140628580
$ws.Range("J86").Value = 27785278
$ws.Range("K86").Value = 140628580
$ws.Range("L86").Value = 27785278
$ws.Range("M86").Value = -140627457
$ws.Range("N86").Value = -27787524
$ws.Range("H89").Value = 118059910
$ws.Range("I89").Value = 140628580
$ws.Range("J89").Value = 27785278
$ws.Range("K89").Value = 703142900
$ws.Range("L89").Value = 138926390
$ws.Range("M89").Value = -703137284
$ws.Range("N89").Value = -138937622
$ws.Range("H132").Value = 1353.7727
$ws.Range("I132").Value = 1322.8572
$ws.Range("K132").Value = 3968.5716
$ws.Range("M132").Value = -1438.5716
$ws.Range("H135").Value = 625964.6
$ws.Range("I135").Value = 715142.9
$ws.Range("K135").Value = 6436286.100000001
$ws.Range("M135").Value = -6433751.100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6450.75
$ws.Range("I2").Value = 1611
$ws.Range("J2").Value = 7142.143
$ws.Range("K2").Value = 1611
$ws.Range("L2").Value = 7142.143
$ws.Range("M2").Value = -1498
$ws.Range("N2").Value = -7368.143
$ws.Range("H32").Value = 1957247.2
$ws.Range("I32").Value = 1988304.2
$ws.Range("K32").Value = 1988304.2
$ws.Range("M32").Value = -1988017.2
$ws.Range("H45").Value = 7118.4
$ws.Range("I45").Value = 2948.2
$ws.Range("J45").Value = 11288.6
$ws.Range("K45").Value = 2948.2
$ws.Range("L45").Value = 11288.6
$ws.Range("M45").Value = -2571.2
$ws.Range("N45").Value = -12042.6
$ws.Range("H61").Value = 5090.304
$ws.Range("I61").Value = 2097.121
$ws.Range("J61").Value = 12688.385
$ws.Range("K61").Value = 2097.121
$ws.Range("L61").Value = 12688.385
$ws.Range("M61").Value = -1885.121
$ws.Range("N61").Value = -13112.385
$ws.Range("H110").Value = 47620532
$ws.Range("I110").Value = 1500
$ws.Range("K110").Value = 1500
$ws.Range("M110").Value = 545
$ws.Range("H116").Value = 6450.75
$ws.Range("I116").Value = 1611
$ws.Range("J116").Value = 7142.143
$ws.Range("K116").Value = 1611
$ws.Range("L116").Value = 7142.143
$ws.Range("M116").Value = 683
$ws.Range("N116").Value = -11730.143
$ws.Range("H132").Value = 8549.679
$ws.Range("I132").Value = 10621.333
$ws.Range("J132").Value = 7568.3687
$ws.Range("K132").Value = 31863.999
$ws.Range("L132").Value = 22705.1061
$ws.Range("M132").Value = -29333.999
$ws.Range("N132").Value = -27765.1061
$ws.Range("H136").Value = 5090.304
$ws.Range("I136").Value = 2097.121
$ws.Range("J136").Value = 12688.385
$ws.Range("K136").Value = 6291.363
$ws.Range("L136").Value = 38065.155
$ws.Range("M136").Value = -3741.363
$ws.Range("N136").Value = -43165.155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6450.75
$ws.Range("I3").Value = 1611
$ws.Range("J3").Value = 7142.143
$ws.Range("K3").Value = 1611
$ws.Range("L3").Value = 7142.143
$ws.Range("M3").Value = -1497
$ws.Range("N3").Value = -7370.143
$ws.Range("H80").Value = 309.2
$ws.Range("I80").Value = 376.33334
$ws.Range("J80").Value = 264.44446
$ws.Range("K80").Value = 376.33334
$ws.Range("L80").Value = 264.44446
$ws.Range("M80").Value = 621.66666
$ws.Range("N80").Value = -2260.44446
$ws.Range("H83").Value = 309.2
$ws.Range("I83").Value = 376.33334
$ws.Range("J83").Value = 264.44446
$ws.Range("K83").Value = 1881.6667
$ws.Range("L83").Value = 1322.2223
$ws.Range("M83").Value = 3110.3333
$ws.Range("N83").Value = -11306.2223
$ws.Range("H98").Value = 53942
$ws.Range("J98").Value = 53942
$ws.Range("L98").Value = 53942
$ws.Range("N98").Value = -59932
$ws.Range("H107").Value = 59213570
$ws.Range("I107").Value = 66178224
$ws.Range("J107").Value = 14000
$ws.Range("K107").Value = 66178224
$ws.Range("L107").Value = 14000
$ws.Range("M107").Value = -66176304
$ws.Range("N107").Value = -17840
$ws.Range("H134").Value = 6300.8535
$ws.Range("I134").Value = 2529.7368
$ws.Range("K134").Value = 7589.2104
$ws.Range("M134").Value = -5054.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11915196
$ws.Range("I31").Value = 5336.625
$ws.Range("K31").Value = 5336.625
$ws.Range("M31").Value = -5041.625
$ws.Range("H34").Value = 11915196
$ws.Range("I34").Value = 5336.625
$ws.Range("K34").Value = 5336.625
$ws.Range("M34").Value = -5134.625
$ws.Range("H60").Value = 17999.166
$ws.Range("I60").Value = 8250
$ws.Range("K60").Value = 8250
$ws.Range("M60").Value = -7739
$ws.Range("H132").Value = 12909276
$ws.Range("I132").Value = 2340.389
$ws.Range("J132").Value = 30780416
$ws.Range("K132").Value = 7021.167
$ws.Range("L132").Value = 92341248
$ws.Range("M132").Value = -4491.167
$ws.Range("N132").Value = -92346308

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3538.6667
$ws.Range("I3").Value = 3538.6667
$ws.Range("K3").Value = 10616.0001
$ws.Range("M3").Value = -10504.0001
$ws.Range("H86").Value = 580.36365
$ws.Range("I86").Value = 598.3333
$ws.Range("J86").Value = 499.5
$ws.Range("K86").Value = 1794.9999
$ws.Range("L86").Value = 1498.5
$ws.Range("M86").Value = -608.9999
$ws.Range("N86").Value = -3870.5
$ws.Range("H89").Value = 580.36365
$ws.Range("I89").Value = 598.3333
$ws.Range("J89").Value = 499.5
$ws.Range("K89").Value = 5384.9997
$ws.Range("L89").Value = 4495.5
$ws.Range("M89").Value = 543.0002999999997
$ws.Range("N89").Value = -16351.5
$ws.Range("H128").Value = 265000.5
$ws.Range("I128").Value = 265000.5
$ws.Range("K128").Value = 795001.5
$ws.Range("M128").Value = -790021.5
$ws.Range("H132").Value = 9476.385
$ws.Range("J132").Value = 13441.857
$ws.Range("L132").Value = 120976.713
$ws.Range("N132").Value = -126036.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 74997.25
$ws.Range("J52").Value = 89996.336
$ws.Range("L52").Value = 89996.336
$ws.Range("N52").Value = -90514.336
$ws.Range("H58").Value = 75799.8
$ws.Range("J58").Value = 75799.8
$ws.Range("L58").Value = 75799.8
$ws.Range("N58").Value = -76353.8
$ws.Range("H113").Value = 5053.778
$ws.Range("I113").Value = 3030.3
$ws.Range("J113").Value = 7583.125
$ws.Range("K113").Value = 3030.3
$ws.Range("L113").Value = 7583.125
$ws.Range("M113").Value = -860.3000000000002
$ws.Range("N113").Value = -11923.125
$ws.Range("H130").Value = 89494
$ws.Range("J130").Value = 89494
$ws.Range("L130").Value = 89494
$ws.Range("N130").Value = -99534
$ws.Range("H134").Value = 119326
$ws.Range("J134").Value = 119326
$ws.Range("L134").Value = 357978
$ws.Range("N134").Value = -363048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 41667092
$ws.Range("I55").Value = 100000080
$ws.Range("J55").Value = 669.0714
$ws.Range("K55").Value = 100000080
$ws.Range("L55").Value = 669.0714
$ws.Range("M55").Value = -99999907
$ws.Range("N55").Value = -1015.0714
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H136").Value = 11095.73
$ws.Range("J136").Value = 12857.143
$ws.Range("L136").Value = 38571.429
$ws.Range("N136").Value = -43671.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14292504
$ws.Range("I81").Value = 2921.5833
$ws.Range("J81").Value = 100030000
$ws.Range("K81").Value = 5843.1666
$ws.Range("L81").Value = 200060000
$ws.Range("M81").Value = -4782.1666
$ws.Range("N81").Value = -200062122
$ws.Range("H84").Value = 14292504
$ws.Range("I84").Value = 2921.5833
$ws.Range("J84").Value = 100030000
$ws.Range("K84").Value = 29215.833
$ws.Range("L84").Value = 1000300000
$ws.Range("M84").Value = -23911.833
$ws.Range("N84").Value = -1000310608
$ws.Range("H107").Value = 708.41174
$ws.Range("I107").Value = 453
$ws.Range("K107").Value = 1359
$ws.Range("M107").Value = 561
$ws.Range("H122").Value = 367948.72
$ws.Range("J122").Value = 6405
$ws.Range("L122").Value = 19215
$ws.Range("N122").Value = -24115
$ws.Range("H132").Value = 55560156
$ws.Range("J132").Value = 5483.857
$ws.Range("L132").Value = 16451.571
$ws.Range("N132").Value = -21511.571
$ws.Range("H136").Value = 34522948
$ws.Range("I136").Value = 125001320
$ws.Range("K136").Value = 375003960
$ws.Range("M136").Value = -375001410
